$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)

$ws.Range("A13").Value = (Get-Date -Year 2018 -Month 5 -Day 26 -Hour 16 -Minute 0 -Second 0)
$ws.Range("B13").Value = (Get-Date -Year 2018 -Month 5 -Day 26 -Hour 18 -Minute 0 -Second 0)
$ws.Range("C13").Value = "Create code to generate final_predict.csv containing ordering in test_set"

$ws.Range("A14").Select()
